$wb = $excel.ActiveWorkbook

# --- Sheet "Produtos": fill in Quantidade / Vendidos / Preço for the two products ---
$wsProdutos = $wb.Worksheets.Item("Produtos")

# Row 2 - "Bacon T"
$wsProdutos.Range("B2").Value = 900
$wsProdutos.Range("C2").Value = 1200
$wsProdutos.Range("D2").Value = 5

# Row 3 - "Cama X"
$wsProdutos.Range("B3").Value = 20
$wsProdutos.Range("C3").Value = 40
$wsProdutos.Range("D3").Value = 550.5

# --- Sheet "Ofertas": add first offer row ---
$wsOfertas = $wb.Worksheets.Item("Ofertas")

$wsOfertas.Range("A2").Value = "Oferta 1"
$wsOfertas.Range("B2").Value = "2015-08-11 00:00:00"
$wsOfertas.Range("C2").Value = "2019-11-11 00:00:00"
